$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.41
$ws.Range("D5").Value = 0.519
$ws.Range("E5").Value = 0.552
$ws.Range("F5").Value = 0.609
$ws.Range("G5").Value = 0.621
$ws.Range("H5").Value = 0.637

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.41
$ws.Range("E7").Value = 0.552
$ws.Range("F7").Value = 0.609
$ws.Range("H7").Value = 0.637

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.398
$ws.Range("D8").Value = 0.589
$ws.Range("E8").Value = 0.625
$ws.Range("F8").Value = 0.662
$ws.Range("G8").Value = 0.6909999999999999
$ws.Range("H8").Value = 0.703

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.305
$ws.Range("C9").Value = 0.474
$ws.Range("D9").Value = 0.602
$ws.Range("E9").Value = 0.628
$ws.Range("F9").Value = 0.642
$ws.Range("G9").Value = 0.672
$ws.Range("H9").Value = 0.684
